$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data (order matters for shared-string table indices):
# C6 (popAbwesenheitAnlegen selector) is inserted first -> shared string index 14
# D4 (dlgPasswordReset title) is inserted second -> shared string index 15
$ws.Range("C6").Value = "//div[@id='createAbsenceModal']"
$ws.Range("D4").Value = "Passwort vergessen | TT-Planer"

# Widen columns C and D so the new, longer values fit
# (mirrors Excel's own "AutoFit Column Width" after the new values were entered)
$ws.Range("C1").ColumnWidth = 28
$ws.Range("D1").ColumnWidth = 26

# Update the active selection to F11
$ws.Range("F11").Select()

# Update window position
$excel.ActiveWindow.Left = 252
$excel.ActiveWindow.Top = 1380
